$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap contents of rows 3 and 4 for the columns that actually differ
# Row 3 <- old Row 4 values ; Row 4 <- old Row 3 values
$row3_A = $ws.Range("A3").Value2
$row3_B = $ws.Range("B3").Value2
$row3_E = $ws.Range("E3").Value2
$row3_F = $ws.Range("F3").Value2
$row3_G = $ws.Range("G3").Value2
$row3_H = $ws.Range("H3").Value2
$row3_Q = $ws.Range("Q3").Value2
$row3_R = $ws.Range("R3").Value2
$row3_Z = $ws.Range("Z3").Value2
$row3_AB = $ws.Range("AB3").Value2

$row4_A = $ws.Range("A4").Value2
$row4_B = $ws.Range("B4").Value2
$row4_E = $ws.Range("E4").Value2
$row4_F = $ws.Range("F4").Value2
$row4_G = $ws.Range("G4").Value2
$row4_H = $ws.Range("H4").Value2
$row4_Q = $ws.Range("Q4").Value2
$row4_R = $ws.Range("R4").Value2
$row4_Z = $ws.Range("Z4").Value2
$row4_AB = $ws.Range("AB4").Value2

$ws.Range("A3").Value = $row4_A
$ws.Range("B3").Value = $row4_B
$ws.Range("E3").Value = $row4_E
$ws.Range("F3").Value = $row4_F
$ws.Range("G3").Value = $row4_G
$ws.Range("H3").Value = $row4_H
$ws.Range("Q3").Value = $row4_Q
$ws.Range("R3").Value = $row4_R
$ws.Range("Z3").Value = $row4_Z
$ws.Range("AB3").Value = $row4_AB

$ws.Range("A4").Value = $row3_A
$ws.Range("B4").Value = $row3_B
$ws.Range("E4").Value = $row3_E
$ws.Range("F4").Value = $row3_F
$ws.Range("G4").Value = $row3_G
$ws.Range("H4").Value = $row3_H
$ws.Range("Q4").Value = $row3_Q
$ws.Range("R4").Value = $row3_R
$ws.Range("Z4").Value = $row3_Z
$ws.Range("AB4").Value = $row3_AB

# Swap contents of rows 14 and 15 for the columns that actually differ
$row14_A = $ws.Range("A14").Value2
$row14_Q = $ws.Range("Q14").Value2
$row14_R = $ws.Range("R14").Value2
$row14_Z = $ws.Range("Z14").Value2
$row14_AB = $ws.Range("AB14").Value2

$row15_A = $ws.Range("A15").Value2
$row15_Q = $ws.Range("Q15").Value2
$row15_R = $ws.Range("R15").Value2
$row15_Z = $ws.Range("Z15").Value2
$row15_AB = $ws.Range("AB15").Value2

$ws.Range("A14").Value = $row15_A
$ws.Range("Q14").Value = $row15_Q
$ws.Range("R14").Value = $row15_R
$ws.Range("Z14").Value = $row15_Z
$ws.Range("AB14").Value = $row15_AB

$ws.Range("A15").Value = $row14_A
$ws.Range("Q15").Value = $row14_Q
$ws.Range("R15").Value = $row14_R
$ws.Range("Z15").Value = $row14_Z
$ws.Range("AB15").Value = $row14_AB
